$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the crypto price/volume refresh.
# Column D ("Price") values are forced to text (quote-prefixed) and the
# cell style is restored afterwards so Excel does not reinterpret values
# such as "1.00" or "577.54" as numbers.
$updates = @(
    @{ Ref = 'D2'; Value = '67.106.01' },
    @{ Ref = 'E2'; Value = '  +4.44%  ' },
    @{ Ref = 'D3'; Value = '3.246.42' },
    @{ Ref = 'E3'; Value = '  +2.16%  ' },
    @{ Ref = 'E4'; Value = '  +0.05%  ' },
    @{ Ref = 'D5'; Value = '577.54' },
    @{ Ref = 'E5'; Value = '  +3.68%  ' },
    @{ Ref = 'D6'; Value = '176.34' },
    @{ Ref = 'E6'; Value = '  +3.26%  ' },
    @{ Ref = 'B7'; Value = 'USDC' },
    @{ Ref = 'C7'; Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc' },
    @{ Ref = 'D7'; Value = '1.00' },
    @{ Ref = 'E7'; Value = '  +0.04%  ' },
    @{ Ref = 'B8'; Value = 'XRP' },
    @{ Ref = 'C8'; Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp' },
    @{ Ref = 'D8'; Value = '0.605' },
    @{ Ref = 'E8'; Value = '  +0.72%  ' },
    @{ Ref = 'D9'; Value = '3.247.80' },
    @{ Ref = 'E9'; Value = '  +2.27%  ' },
    @{ Ref = 'E10'; Value = '  +4.73%  ' },
    @{ Ref = 'E11'; Value = '  +1.45%  ' },
    @{ Ref = 'D12'; Value = '0.409' },
    @{ Ref = 'E12'; Value = '  +3.32%  ' },
    @{ Ref = 'D13'; Value = '3.812.06' },
    @{ Ref = 'E13'; Value = '  +2.47%  ' },
    @{ Ref = 'E14'; Value = '  +1.90%  ' },
    @{ Ref = 'D15'; Value = '27.83' },
    @{ Ref = 'E15'; Value = '  +1.20%  ' },
    @{ Ref = 'D16'; Value = '67.088.29' },
    @{ Ref = 'E16'; Value = '  +4.33%  ' },
    @{ Ref = 'E17'; Value = '  +3.86%  ' },
    @{ Ref = 'D18'; Value = '3.248.44' },
    @{ Ref = 'E18'; Value = '  +2.87%  ' },
    @{ Ref = 'D19'; Value = '5.81' },
    @{ Ref = 'E19'; Value = '  +3.03%  ' },
    @{ Ref = 'D20'; Value = '13.29' },
    @{ Ref = 'E20'; Value = '  +2.44%  ' },
    @{ Ref = 'D21'; Value = '367.63' },
    @{ Ref = 'E21'; Value = '  +4.33%  ' },
    @{ Ref = 'D22'; Value = '7.46' },
    @{ Ref = 'E22'; Value = '  +3.81%  ' },
    @{ Ref = 'E23'; Value = '  +0.21%  ' },
    @{ Ref = 'D24'; Value = '70.12' },
    @{ Ref = 'E24'; Value = '  +1.61%  ' },
    @{ Ref = 'B25'; Value = 'Polygon' },
    @{ Ref = 'C25'; Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic' },
    @{ Ref = 'D25'; Value = '0.506' },
    @{ Ref = 'E25'; Value = '  +0.97%  ' },
    @{ Ref = 'B26'; Value = 'WrappedeETH' },
    @{ Ref = 'C26'; Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth' },
    @{ Ref = 'D26'; Value = '3.382.55' },
    @{ Ref = 'E26'; Value = '  +2.03%  ' },
    @{ Ref = 'B27'; Value = 'PEPE' },
    @{ Ref = 'C27'; Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe' },
    @{ Ref = 'D27'; Value = '0.0000119' },
    @{ Ref = 'E27'; Value = '  +1.10%  ' },
    @{ Ref = 'D28'; Value = '9.75' },
    @{ Ref = 'E28'; Value = '  +3.97%  ' },
    @{ Ref = 'E29'; Value = '  +1.40%  ' },
    @{ Ref = 'E30'; Value = '  +0.42%  ' },
    @{ Ref = 'E31'; Value = '  +4.05%  ' },
    @{ Ref = 'E32'; Value = '  +0.92%  ' },
    @{ Ref = 'D33'; Value = '22.46' },
    @{ Ref = 'E33'; Value = '  +1.70%  ' },
    @{ Ref = 'E34'; Value = '  -0.22%  ' },
    @{ Ref = 'E35'; Value = '  +3.95%  ' },
    @{ Ref = 'D36'; Value = '6.77' },
    @{ Ref = 'E36'; Value = '  +2.88%  ' },
    @{ Ref = 'D37'; Value = '171.23' },
    @{ Ref = 'E37'; Value = '  +8.27%  ' },
    @{ Ref = 'D38'; Value = '1.51' },
    @{ Ref = 'E38'; Value = '  +5.49%  ' },
    @{ Ref = 'D39'; Value = '0.850' },
    @{ Ref = 'E39'; Value = '  +5.61%  ' },
    @{ Ref = 'E40'; Value = '  +11.20%  ' },
    @{ Ref = 'D41'; Value = '26.66' },
    @{ Ref = 'E41'; Value = '  +2.39%  ' },
    @{ Ref = 'E42'; Value = '  +4.18%  ' },
    @{ Ref = 'B43'; Value = 'Maker' },
    @{ Ref = 'C43'; Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr' },
    @{ Ref = 'D43'; Value = '2.703.49' },
    @{ Ref = 'E43'; Value = '  +1.86%  ' },
    @{ Ref = 'B44'; Value = 'RenderToken' },
    @{ Ref = 'C44'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' },
    @{ Ref = 'D44'; Value = '6.33' },
    @{ Ref = 'E44'; Value = '  +5.49%  ' },
    @{ Ref = 'E45'; Value = '  +3.82%  ' },
    @{ Ref = 'D46'; Value = '40.39' },
    @{ Ref = 'E46'; Value = '  +4.51%  ' },
    @{ Ref = 'D47'; Value = '0.0671' },
    @{ Ref = 'E47'; Value = '  +2.77%  ' },
    @{ Ref = 'D48'; Value = '24.62' },
    @{ Ref = 'E48'; Value = '  +5.69%  ' },
    @{ Ref = 'D49'; Value = '333.10' },
    @{ Ref = 'E49'; Value = '  +3.63%  ' },
    @{ Ref = 'E50'; Value = '  +3.82%  ' },
    @{ Ref = 'E51'; Value = '  +2.21%  ' }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Ref)
    if ($u.Ref.Substring(0,1) -eq 'D') {
        $origStyle = $range.Style
        $range.Value = "'" + $u.Value
        $range.Style = $origStyle
    } else {
        $range.Value = $u.Value
    }
}
